$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated Volue Forecast (MW) values in column D for all data rows,
# and fill in previously-missing Hydro Reservoir/River/Actual values for rows 46-48
# (api_key_entsoe missing value fix).

$ws.Range("D2").Value = 1738.388764554137
$ws.Range("D3").Value = 1727.258761081866
$ws.Range("D4").Value = 1716.128757609596
$ws.Range("D5").Value = 1704.998754137325
$ws.Range("D6").Value = 1703.656087051781
$ws.Range("D7").Value = 1702.313419966238
$ws.Range("D8").Value = 1700.970752880694
$ws.Range("D9").Value = 1699.62808579515
$ws.Range("D10").Value = 1701.076752913815
$ws.Range("D11").Value = 1702.525420032479
$ws.Range("D12").Value = 1703.974087151143
$ws.Range("D13").Value = 1705.422754269807
$ws.Range("D14").Value = 1716.482091053419
$ws.Range("D15").Value = 1727.54142783703
$ws.Range("D16").Value = 1738.600764620641
$ws.Range("D17").Value = 1749.660101404253
$ws.Range("D18").Value = 1762.874772193393
$ws.Range("D19").Value = 1776.089442982533
$ws.Range("D20").Value = 1789.304113771673
$ws.Range("D21").Value = 1802.518784560813
$ws.Range("D22").Value = 1867.273026984722
$ws.Range("D23").Value = 1932.027269408632
$ws.Range("D24").Value = 1996.781511832542
$ws.Range("D25").Value = 2061.535754256451
$ws.Range("D26").Value = 2050.523528598762
$ws.Range("D27").Value = 2039.511302941074
$ws.Range("D28").Value = 2028.499077283385
$ws.Range("D29").Value = 2017.486851625697
$ws.Range("D30").Value = 1961.224389628772
$ws.Range("D31").Value = 1904.961927631848
$ws.Range("D32").Value = 1848.699465634923
$ws.Range("D33").Value = 1792.437003637999
$ws.Range("D34").Value = 1757.327437129125
$ws.Range("D35").Value = 1722.217870620252
$ws.Range("D36").Value = 1687.108304111379
$ws.Range("D37").Value = 1651.998737602505
$ws.Range("D38").Value = 1609.634057719259
$ws.Range("D39").Value = 1567.269377836012
$ws.Range("D40").Value = 1524.904697952765
$ws.Range("D41").Value = 1482.540018069518
$ws.Range("D42").Value = 1485.260685584962
$ws.Range("D43").Value = 1487.981353100407
$ws.Range("D44").Value = 1490.702020615851
$ws.Range("D45").Value = 1493.422688131295
$ws.Range("B46").Value = 106
$ws.Range("C46").Value = 1129
$ws.Range("D46").Value = 1492.621798992549
$ws.Range("E46").Value = 1235
$ws.Range("B47").Value = 106
$ws.Range("C47").Value = 1115
$ws.Range("D47").Value = 1491.820909853804
$ws.Range("E47").Value = 1221
$ws.Range("B48").Value = 107
$ws.Range("C48").Value = 1124
$ws.Range("D48").Value = 1491.020020715059
$ws.Range("E48").Value = 1231
$ws.Range("D49").Value = 1490.219131576313
$ws.Range("D50").Value = 1497.745133924139
$ws.Range("D51").Value = 1505.271136271965
$ws.Range("D52").Value = 1512.797138619791
$ws.Range("D53").Value = 1520.323140967617
$ws.Range("D54").Value = 1541.181591919447
$ws.Range("D55").Value = 1562.040042871278
$ws.Range("D56").Value = 1582.898493823108
$ws.Range("D57").Value = 1603.756944774938
$ws.Range("D58").Value = 1648.241625319616
$ws.Range("D59").Value = 1692.726305864293
$ws.Range("D60").Value = 1737.210986408971
$ws.Range("D61").Value = 1781.695666953649
$ws.Range("D62").Value = 1866.213026654256
$ws.Range("D63").Value = 1950.730386354863
$ws.Range("D64").Value = 2035.24774605547
$ws.Range("D65").Value = 2119.765105756077
$ws.Range("D66").Value = 2179.619791095779
$ws.Range("D67").Value = 2239.474476435481
$ws.Range("D68").Value = 2299.329161775183
$ws.Range("D69").Value = 2359.183847114886
$ws.Range("D70").Value = 2399.781859780482
$ws.Range("D71").Value = 2440.379872446078
$ws.Range("D72").Value = 2480.977885111674
$ws.Range("D73").Value = 2521.57589777727
$ws.Range("D74").Value = 2522.553453637797
$ws.Range("D75").Value = 2523.531009498325
$ws.Range("D76").Value = 2524.508565358852
$ws.Range("D77").Value = 2525.48612121938
$ws.Range("D78").Value = 2500.446557852131
$ws.Range("D79").Value = 2475.406994484882
$ws.Range("D80").Value = 2450.367431117634
$ws.Range("D81").Value = 2425.327867750385
$ws.Range("D82").Value = 2399.746526436338
$ws.Range("D83").Value = 2374.16518512229
$ws.Range("D84").Value = 2348.583843808244
$ws.Range("D85").Value = 2323.002502494197
$ws.Range("D86").Value = 2245.0218114994
$ws.Range("D87").Value = 2167.041120504603
$ws.Range("D88").Value = 2089.060429509806
$ws.Range("D89").Value = 2011.079738515009
$ws.Range("D90").Value = 1975.298838463702
$ws.Range("D91").Value = 1939.517938412395
$ws.Range("D92").Value = 1903.737038361088
$ws.Range("D93").Value = 1867.956138309781
$ws.Range("D94").Value = 1858.776954463491
$ws.Range("D95").Value = 1849.597770617201
$ws.Range("D96").Value = 1840.418586770911
$ws.Range("D97").Value = 1831.239402924621
